# Edit: rename "总计" sheet to "2022-Q1" (repurposed with fund-level holdings
# data), and add a brand-new "总计" sheet after it carrying the updated
# date/count/value summary table (with a new 2022-Q1 row inserted on top).

$wb = $excel.ActiveWorkbook

# ---- locate existing sheets --------------------------------------------
$sheetCount = $wb.Worksheets.Count
$oldTotal = $wb.Worksheets.Item($sheetCount)          # current "总计" sheet
$lastQuarter = $wb.Worksheets.Item($sheetCount - 1)    # current "2021-Q4" sheet

# ---- 1) create the new "总计" sheet as a copy of the old one BEFORE we
#         touch anything, so it inherits all header/row formatting -------
$oldTotal.Copy($null, $oldTotal)
$newTotal = $wb.Worksheets.Item($sheetCount + 1)

# ---- 2) rename the original sheet to "2022-Q1" (freeing up the "总计"
#         name) and rebuild its content with the per-fund holdings table,
#         then rename the copy to "总计" ----------------------------------
$oldTotal.Name = "2022-Q1"
$q1 = $oldTotal
$newTotal.Name = "总计"

# Bring in the exact column layout/styling used by the other quarterly
# sheets (header + data-row formatting) by copying the 2021-Q4 sheet's
# whole used range over first.
$lastQuarter.Range("A1:H34").Copy($q1.Range("A1"))

$fundRows = @(
    @("004224", "南方军工改革灵活配置混合A", "63.35", "86.90", "7.87", "4.9856", 8),
    @("002251", "华夏军工安全灵活配置混合", "44.78", "94.71", "9.20", "4.1198", 1),
    @("000001", "华夏成长混合", "31.69", "72.21", "5.25", "1.6637", 2),
    @("002983", "长信国防军工量化灵活配置混合A", "26.74", "93.30", "4.96", "1.3263", 10),
    @("011148", "南方军工改革灵活配置混合C", "15.56", "86.90", "7.87", "1.2246", 8),
    @("002345", "华夏高端制造灵活配置混合", "23.47", "93.20", "4.16", "0.9764", 6),
    @("010410", "长城品质成长混合A", "29.08", "70.65", "2.41", "0.7008", 5),
    @("008960", "长信国防军工量化灵活配置混合C", "13.19", "93.30", "4.96", "0.6542", 10),
    @("580009", "东吴多策略灵活配置混合", "8.31", "88.46", "5.74", "0.4770", 6),
    @("200007", "长城安心回报混合", "11.42", "71.80", "3.93", "0.4488", 2),
    @("200012", "长城中小盘成长混合", "12.65", "84.26", "3.13", "0.3959", 3),
    @("011455", "长城竞争优势六个月持有期混合型证券投资基金A", "9.99", "79.15", "3.90", "0.3896", 4),
    @("160143", "南方创业板2年定期开放混合", "8.23", "84.71", "4.41", "0.3629", 5),
    @("000977", "长城环保主题灵活配置混合", "9.63", "81.27", "3.10", "0.2985", 7),
    @("005738", "长城智能产业灵活配置混合", "8.04", "84.35", "3.44", "0.2766", 3),
    @("010049", "长城成长先锋混合A", "7.27", "81.13", "3.09", "0.2246", 7),
    @("006769", "长城研究精选混合", "6.41", "82.49", "2.55", "0.1635", 4),
    @("014189", "南方专精特新混合A", "4.99", "30.05", "2.89", "0.1442", 2),
    @("002703", "长城久源灵活配置混合", "1.56", "89.80", "7.95", "0.1240", 1),
    @("002542", "长城久鼎灵活配置混合", "3.81", "81.60", "3.13", "0.1193", 6),
    @("010411", "长城品质成长混合C", "2.34", "70.65", "2.41", "0.0564", 5),
    @("002885", "摩根士丹利华鑫万众创新灵活配置混合", "0.86", "94.36", "5.11", "0.0439", 7),
    @("014190", "南方专精特新混合C", "1.13", "30.05", "2.89", "0.0327", 2),
    @("004091", "博时沪港深价值优选灵活配置混合A", "1.28", "58.32", "2.38", "0.0305", 3),
    @("010050", "长城成长先锋混合C", "0.90", "81.13", "3.09", "0.0278", 7),
    @("011456", "长城竞争优势六个月持有期混合型证券投资基金C", "0.65", "79.15", "3.90", "0.0254", 4),
    @("002145", "诺安景鑫灵活配置混合", "0.53", "83.45", "3.65", "0.0193", 10),
    @("501002", "长信价值优选混合", "0.46", "81.36", "1.89", "0.0087", 5),
    @("162107", "金鹰量化精选股票（LOF）", "0.06", "93.77", "5.89", "0.0035", 6),
    @("009327", "东兴兴晟混合A", "0.34", "79.83", "0.92", "0.0031", 7),
    @("004092", "博时沪港深价值优选灵活配置混合C", "0.09", "58.32", "2.38", "0.0021", 3),
    @("009328", "东兴兴晟混合C", "0.08", "79.83", "0.92", "0.0007", 7),
    @("006992", "嘉合锦创优势精选混合", "0.02", "74.79", "2.33", "0.0005", 8)
)

$r = 2
foreach ($row in $fundRows) {
    $code = $row[0]
    $name = $row[1]
    $scale = $row[2]
    $totalPos = $row[3]
    $posRatio = $row[4]
    $heldValue = $row[5]
    $rank = $row[6]

    $q1.Cells.Item($r, 1).Value = $r - 2
    $q1.Cells.Item($r, 2).Value = "'" + $code
    $q1.Cells.Item($r, 3).Value = $name
    $q1.Cells.Item($r, 4).Value = "'" + $scale
    $q1.Cells.Item($r, 5).Value = "'" + $totalPos
    $q1.Cells.Item($r, 6).Value = "'" + $posRatio
    $q1.Cells.Item($r, 7).Value = "'" + $heldValue
    $q1.Cells.Item($r, 8).Value = $rank

    # The apostrophe-prefixed assignments above force genuine text storage
    # (preserving leading zeros / trailing zeros), but Excel also stamps a
    # quote-prefixed "text" number format on those cells. Re-pull the
    # formatting from the fund-name cell (column C), which is plain,
    # unstyled text, so the numeric-looking text cells end up with no
    # special styling -- matching the other quarterly sheets.
    $q1.Cells.Item($r, 3).Copy()
    $q1.Cells.Item($r, 2).PasteSpecial(-4122)
    $q1.Cells.Item($r, 4).PasteSpecial(-4122)
    $q1.Cells.Item($r, 5).PasteSpecial(-4122)
    $q1.Cells.Item($r, 6).PasteSpecial(-4122)
    $q1.Cells.Item($r, 7).PasteSpecial(-4122)

    $r = $r + 1
}
$excel.CutCopyMode = $false

# ---- 3) rebuild the new "总计" sheet: insert a 2022-Q1 row on top and
#         shift the existing summary rows down ---------------------------
$newTotal.Rows.Item(2).Insert()
$newTotal.Range("A2:D2").ClearFormats()

# Restore A2's header-index styling (bold/centered) by pulling the format
# from A3 (the row pushed down, which already carries that style).
$newTotal.Cells.Item(3, 1).Copy()
$newTotal.Cells.Item(2, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false

$newTotal.Cells.Item(2, 1).Value = 0
$newTotal.Cells.Item(2, 2).Value = "2022-Q1"
$newTotal.Cells.Item(2, 3).Value = 33
$newTotal.Cells.Item(2, 4).Value = 19.33

# Renumber the index column for the rows that shifted down.
$newTotal.Cells.Item(3, 1).Value = 1
$newTotal.Cells.Item(4, 1).Value = 2
$newTotal.Cells.Item(5, 1).Value = 3
$newTotal.Cells.Item(6, 1).Value = 4
$newTotal.Cells.Item(7, 1).Value = 5
